# The reviewer annotated (side-by-side score, column E) rows 33 through 55,
# extending the previously-annotated range (rows 2-32) all the way to the
# bottom of the sheet (row 55).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> annotation score (-2 .. 2)
$scores = [ordered]@{
    33 = 0
    34 = -1
    35 = 0
    36 = -2
    37 = -2
    38 = 0
    39 = 0
    40 = 0
    41 = 0
    42 = 1
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = -1
    49 = -2
    50 = 1
    51 = 0
    52 = 0
    53 = 0
    54 = 0
    55 = 0
}

foreach ($row in $scores.Keys) {
    $ws.Cells.Item($row, 5).Value = $scores[$row]
}

# Leave the selection on the last cell that was filled in, matching the
# scrolled-down state the sheet was left in after annotating through row 55.
$ws.Range("E55").Select()
